$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Periodo Mora updated from 2507 to 2508 for every worker row (E16:E21)
$ws.Range("E16:E21").Value = "2508"

# 2) Row 19 previously held PPT / 1999120 / ANTONY MOISES RINCON MATOS.
#    That worker was replaced; the row now mirrors row 18 (CC / 1148147776 / YAQUELIN OSPINO ACUÑA).
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1148147776"
$ws.Range("D19").Value = "YAQUELIN OSPINO ACUÑA"
